# Atualizei dados da bibi
# - corrige o valor de total_venda do dia 5 (06/2025)
# - adiciona um novo dia (dia 10) ao bloco de 06/2025, empurrando as linhas seguintes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige o valor de faturamento do dia 5 (06/2025), linha 6
$ws.Range("B6").Value = 21643.35

# Insere uma nova linha na posicao 11 (antes do bloco de 05/2025),
# deslocando as linhas existentes uma posicao para baixo
$ws.Rows.Item(11).Insert()

# Preenche a nova linha 11 com o dia 10 do periodo 06/2025
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 31200
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 2025
$ws.Range("E11").Value = "06/2025"
